$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '41.631.30'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '2.463.92'
$ws.Range('E4').Value = '  +0.52%  '
$r = $ws.Range('D5')
$r.Value = "'314.84"
$r.Style = "Normal"
$ws.Range('E5').Value = '  +0.49%  '
$r = $ws.Range('D6')
$r.Value = "'92.24"
$r.Style = "Normal"
$ws.Range('E6').Value = '  -2.46%  '
$r = $ws.Range('D7')
$r.Value = "'0.548"
$r.Style = "Normal"
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  +0.49%  '
$ws.Range('E9').Value = '  +2.43%  '
$r = $ws.Range('D10')
$r.Value = "'32.43"
$r.Style = "Normal"
$ws.Range('E10').Value = '  -1.62%  '
$r = $ws.Range('D11')
$r.Value = "'0.0834"
$r.Style = "Normal"
$ws.Range('E11').Value = '  +6.05%  '
$ws.Range('E12').Value = '  +0.38%  '
$ws.Range('D13').Value = '2.841.37'
$ws.Range('E13').Value = '  -1.22%  '
$r = $ws.Range('D14')
$r.Value = "'6.86"
$r.Style = "Normal"
$ws.Range('E14').Value = '  -0.28%  '
$r = $ws.Range('D15')
$r.Value = "'15.86"
$r.Style = "Normal"
$ws.Range('E15').Value = '  +2.12%  '
$ws.Range('D16').Value = '2.452.60'
$ws.Range('E16').Value = '  -3.53%  '
$ws.Range('E17').Value = '  +2.25%  '
$ws.Range('D18').Value = '41.598.59'
$ws.Range('E18').Value = '  -0.17%  '
$r = $ws.Range('D19')
$r.Value = "'6.48"
$r.Style = "Normal"
$ws.Range('E19').Value = '  +2.07%  '
$ws.Range('D20').Value = '0.0₃0949'
$ws.Range('E20').Value = '  +2.63%  '
$r = $ws.Range('D21')
$r.Value = "'70.66"
$r.Style = "Normal"
$ws.Range('E21').Value = '  +0.02%  '
$r = $ws.Range('D22')
$r.Value = "'11.33"
$r.Style = "Normal"
$ws.Range('E22').Value = '  +0.99%  '
$r = $ws.Range('D23')
$r.Value = "'238.72"
$r.Style = "Normal"
$ws.Range('E23').Value = '  +1.03%  '
$ws.Range('E24').Value = '  -0.21%  '
$ws.Range('E25').Value = '  +0.41%  '
$ws.Range('E26').Value = '  -0.08%  '
$r = $ws.Range('D27')
$r.Value = "'24.39"
$r.Style = "Normal"
$ws.Range('E27').Value = '  -1.16%  '
$r = $ws.Range('D28')
$r.Value = "'2.25"
$r.Style = "Normal"
$ws.Range('E28').Value = '  +0.40%  '
$r = $ws.Range('D29')
$r.Value = "'9.73"
$r.Style = "Normal"
$ws.Range('E29').Value = '  +0.35%  '
$r = $ws.Range('D30')
$r.Value = "'35.16"
$r.Style = "Normal"
$ws.Range('E30').Value = '  -3.17%  '
$ws.Range('E31').Value = '  +0.80%  '
$r = $ws.Range('D32')
$r.Value = "'5.47"
$r.Style = "Normal"
$ws.Range('E32').Value = '  +0.62%  '
$r = $ws.Range('D33')
$r.Value = "'2.57"
$r.Style = "Normal"
$ws.Range('E33').Value = '  -0.10%  '
$r = $ws.Range('D34')
$r.Value = "'0.0760"
$r.Style = "Normal"
$ws.Range('E34').Value = '  +0.21%  '
$ws.Range('E35').Value = '  -0.90%  '
$ws.Range('E36').Value = '  -5.11%  '
$ws.Range('E37').Value = '  -3.47%  '
$ws.Range('E38').Value = '  +1.02%  '
$ws.Range('E39').Value = '  +1.84%  '
$ws.Range('E40').Value = '  -2.50%  '
$r = $ws.Range('D41')
$r.Value = "'3.96"
$r.Style = "Normal"
$ws.Range('E41').Value = '  -4.62%  '
$ws.Range('E42').Value = '  +0.44%  '
$ws.Range('D43').Value = '1.975.55'
$ws.Range('E43').Value = '  +0.96%  '
$ws.Range('B44').Value = 'VeChain'
$ws.Range('C44').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$r = $ws.Range('D44')
$r.Value = "'0.0282"
$r.Style = "Normal"
$ws.Range('E44').Value = '  -0.92%  '
$ws.Range('B45').Value = 'EnergySwap'
$ws.Range('C45').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$r = $ws.Range('D45')
$r.Value = "'18.92"
$r.Style = "Normal"
$ws.Range('E45').Value = '  -4.36%  '
$ws.Range('E46').Value = '  -2.07%  '
$r = $ws.Range('D47')
$r.Value = "'9.03"
$r.Style = "Normal"
$ws.Range('E47').Value = '  +2.12%  '
$ws.Range('D48').Value = '2.697.71'
$ws.Range('E48').Value = '  -1.32%  '
$r = $ws.Range('D49')
$r.Value = "'96.93"
$r.Style = "Normal"
$ws.Range('E49').Value = '  +0.23%  '
$r = $ws.Range('D50')
$r.Value = "'66.76"
$r.Style = "Normal"
$ws.Range('E50').Value = '  -1.11%  '
$r = $ws.Range('D51')
$r.Value = "'52.45"
$r.Style = "Normal"
$ws.Range('E51').Value = '  +2.92%  '
